# "ranking was trouble shooted"
# Insert a new "Ziel Stadtteil" column (C) between the existing "Stadtteil"
# (A) and "Anzahl Ziele" (old C, becomes D) columns. The new column mirrors
# the neighborhood name already present in column A for each data row.
# Also rename the "Stadtteil" header to "Start Stadtteil" and extend the
# autofilter / filter-database range from A1:C28 to A1:D28.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 28

# Insert a new blank column before the current column C ("Anzahl Ziele"),
# shifting it to column D. Excel carries the column's formatting/width
# along with it, so D keeps its original width of 14.
$ws.Columns.Item(3).Insert()

# New column C ("Ziel Stadtteil") gets the same width as column A (24).
$ws.Columns.Item(3).ColumnWidth = 24 - 5/6

# Update header text.
$ws.Range("A1").Value = "Start Stadtteil"
$ws.Range("C1").Value = "Ziel Stadtteil"

# Populate the new column C with the same Stadtteil name as column A for
# each data row.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = $ws.Cells.Item($r, 1).Value2
}

# Re-apply the autofilter over the new, wider range A1:D28.
$ws.AutoFilterMode = $false
[void]$ws.Range("A1:D28").AutoFilter()

# Keep the hidden _FilterDatabase defined name in sync with the new range.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "='Sheet1'!`$A`$1:`$D`$28"
    }
}
